# Regenerate all penyata to follow new data and format
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Kali X" entries to "Semakan Kali X" for the 3 repeated mini-tables ---
$ws.Range("C16").Value = "Semakan Kali Pertama"
$ws.Range("C17").Value = "Semakan Kali Kedua"
$ws.Range("C18").Value = "Semakan Kali Ketiga"
$ws.Range("C19").Value = "Semakan Kali Keempat"

$ws.Range("C22").Value = "Semakan Kali Pertama"
$ws.Range("C23").Value = "Semakan Kali Kedua"
$ws.Range("C24").Value = "Semakan Kali Ketiga"
$ws.Range("C25").Value = "Semakan Kali Keempat"

$ws.Range("C28").Value = "Semakan Kali Pertama"
$ws.Range("C29").Value = "Semakan Kali Kedua"
$ws.Range("C30").Value = "Semakan Kali Ketiga"
$ws.Range("C31").Value = "Semakan Kali Keempat"

# --- Normalise competition names to title case, and add new competitions ---
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Kad Raya Untuk Guruku"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"
$ws.Range("C38").Value = "Lompat Getah"
$ws.Range("C39").Value = "Theme Party"
$ws.Range("C40").Value = "Hari Koperasi"

# --- Adjust merged cell layout ---
$ws.Range("B15:C15").UnMerge()

# The "STATEMENT OF HOMEROOM ACCOUNT" title moves from E4 into D4 (new merge anchor)
$ws.Range("D4").Value = $ws.Range("E4").Value2
$ws.Range("E4").Value = ""
$ws.Range("D4:G4").Merge()

$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("B43:E43").Merge()

# --- Remove the extra trailing blank row so the sheet ends at row 1000 ---
$ws.Rows.Item(1001).Delete()

# --- Page setup: fit to a single page and centre horizontally when printing ---
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
